# IPWarmupPlan-Sample.xlsx edit:
# Remove the two leading "Properties"/"Value" label rows from the
# "Warmup Plan" sheet (rows 1-2), which were not really part of the
# Phase/Run/Gmail/... data table. Deleting these entire rows shifts the
# rest of the sheet up by two rows (dimension A1:T420 -> A1:T418, merged
# cell ranges shift accordingly) and Excel automatically drops the two
# now-unused shared strings ("Properties" and "Value"), renumbering the
# remaining shared-string indices used on both sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warmup Plan")

# Delete rows 1 and 2 entirely (shifts everything else up by two rows).
$ws.Rows("1:2").Delete()

# Leave the selection on the new header row, matching the state Excel is
# left in right after deleting selected whole rows.
$ws.Rows("1:1").Select()
